$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 37; this shifts the existing rows
# 37..116 down to 38..117 (matches the diff: old row37 -> new row38, ...,
# old row116 -> new row117).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted (blank) row 37 with the new record.
$ws.Cells.Item(37, 1).Value = 6
$ws.Cells.Item(37, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(37, 3).Value = "Metropolitana"
$ws.Cells.Item(37, 4).Value = 44998
$ws.Cells.Item(37, 5).Value = 13
$ws.Cells.Item(37, 6).Value = 100114007
$ws.Cells.Item(37, 7).Value = "Jengibre"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 14000
$ws.Cells.Item(37, 12).Value = 15000
$ws.Cells.Item(37, 13).Value = 14425
$ws.Cells.Item(37, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(37, 15).Value = "Perú"
$ws.Cells.Item(37, 16).Value = 1110
$ws.Cells.Item(37, 17).Value = 13
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest
# of the "Fecha" column.
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
